# Adding the changes we made on may 9th
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows right after the header row (row 1), before the existing
# data (old row 2 onward), pushing old data down by 2 rows.
$ws.Rows.Item(2).Resize(2).Insert()
$ws.Rows.Item(2).Resize(2).ClearFormats()

$newTopRows = @(
    @(0.1050096067542932, -1.756468223065746, 0.4945203567645989),
    @(-0.0286234012063665, -0.7998002785809195, 0.0811297598541999)
)

for ($i = 0; $i -lt $newTopRows.Count; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $newTopRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newTopRows[$i][1]
    $ws.Cells.Item($r, 3).Value = $newTopRows[$i][2]
}

# Append 8 new rows at the end (rows 24-31 after the insert above).
$newBottomRows = @(
    @(-0.4081483519807154, -0.6726997543354425, -0.2190668820118418),
    @(0.2211332225373814, 0.241335413285664, 0.08368853798934378),
    @(0.06768137718341787, 0.3379019900244107, 0.1505034766635118),
    @(0.07254024853511698, 0.5556785336562575, -0.05807583201296457),
    @(0.1816357883567719, 0.1322741392923868, -0.08515337003128903),
    @(-0.02734556931013958, -0.1169588795425942, 0.04497027853313797),
    @(-0.02540700723017953, -0.06986615411481072, -0.074921377335808),
    @(0.02237761537639455, -0.07008743807863513, -0.003453258577050004)
)

$startRow = 24
for ($i = 0; $i -lt $newBottomRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newBottomRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newBottomRows[$i][1]
    $ws.Cells.Item($r, 3).Value = $newBottomRows[$i][2]
}
